$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "61-44=17"
$t.Cell(1,2).Range.Text = "22-6=16"
$t.Cell(1,3).Range.Text = "23-18=5"
$t.Cell(1,4).Range.Text = "6+86=92"
$t.Cell(1,5).Range.Text = "56+35=91"

$t.Cell(2,1).Range.Text = "60-58=2"
$t.Cell(2,2).Range.Text = "90-36=54"
$t.Cell(2,3).Range.Text = "59+25=84"
$t.Cell(2,4).Range.Text = "74-67=7"
$t.Cell(2,5).Range.Text = "3+19=22"

$t.Cell(3,1).Range.Text = "88+8=96"
$t.Cell(3,2).Range.Text = "77-29=48"
$t.Cell(3,3).Range.Text = "8+35=43"
$t.Cell(3,4).Range.Text = "49+18=67"
$t.Cell(3,5).Range.Text = "23+28=51"

$t.Cell(4,1).Range.Text = "90-58=32"
$t.Cell(4,2).Range.Text = "17+58=75"
$t.Cell(4,3).Range.Text = "80-35=45"
$t.Cell(4,4).Range.Text = "44-28=16"
$t.Cell(4,5).Range.Text = "62-44=18"

$t.Cell(5,1).Range.Text = "86-48=38"
$t.Cell(5,2).Range.Text = "46-39=7"
$t.Cell(5,3).Range.Text = "8+45=53"
$t.Cell(5,4).Range.Text = "80-77=3"
$t.Cell(5,5).Range.Text = "49+26=75"

$t.Cell(6,1).Range.Text = "48+5=53"
$t.Cell(6,2).Range.Text = "65+18=83"
$t.Cell(6,3).Range.Text = "65-58=7"
$t.Cell(6,4).Range.Text = "55-6=49"
$t.Cell(6,5).Range.Text = "19+46=65"

$t.Cell(7,1).Range.Text = "55-17=38"
$t.Cell(7,2).Range.Text = "17+17=34"
$t.Cell(7,3).Range.Text = "8+5=13"
$t.Cell(7,4).Range.Text = "10-9=1"
$t.Cell(7,5).Range.Text = "41-2=39"

$t.Cell(8,1).Range.Text = "35+59=94"
$t.Cell(8,2).Range.Text = "36+35=71"
$t.Cell(8,3).Range.Text = "7+24=31"
$t.Cell(8,4).Range.Text = "73-8=65"
$t.Cell(8,5).Range.Text = "39+22=61"

$t.Cell(9,1).Range.Text = "25+36=61"
$t.Cell(9,2).Range.Text = "84-75=9"
$t.Cell(9,3).Range.Text = "76-7=69"
$t.Cell(9,4).Range.Text = "75+17=92"
$t.Cell(9,5).Range.Text = "71-47=24"

$t.Cell(10,1).Range.Text = "76-69=7"
$t.Cell(10,2).Range.Text = "40-16=24"
$t.Cell(10,3).Range.Text = "4+18=22"
$t.Cell(10,4).Range.Text = "47-28=19"
$t.Cell(10,5).Range.Text = "35+19=54"

$t.Cell(11,1).Range.Text = "3+39=42"
$t.Cell(11,2).Range.Text = "64+7=71"
$t.Cell(11,3).Range.Text = "19+3=22"
$t.Cell(11,4).Range.Text = "82-56=26"
$t.Cell(11,5).Range.Text = "35-16=19"

$t.Cell(12,1).Range.Text = "67-19=48"
$t.Cell(12,2).Range.Text = "24+19=43"
$t.Cell(12,3).Range.Text = "19+26=45"
$t.Cell(12,4).Range.Text = "18+39=57"
$t.Cell(12,5).Range.Text = "57+19=76"

$t.Cell(13,1).Range.Text = "69+14=83"
$t.Cell(13,2).Range.Text = "43+48=91"
$t.Cell(13,3).Range.Text = "94-38=56"
$t.Cell(13,4).Range.Text = "59+36=95"
$t.Cell(13,5).Range.Text = "51-29=22"

$t.Cell(14,1).Range.Text = "46-7=39"
$t.Cell(14,2).Range.Text = "54+17=71"
$t.Cell(14,3).Range.Text = "52-16=36"
$t.Cell(14,4).Range.Text = "29+38=67"
$t.Cell(14,5).Range.Text = "54-47=7"

$t.Cell(15,1).Range.Text = "68+3=71"
$t.Cell(15,2).Range.Text = "52-44=8"
$t.Cell(15,3).Range.Text = "61-16=45"
$t.Cell(15,4).Range.Text = "68+29=97"
$t.Cell(15,5).Range.Text = "80-6=74"

$t.Cell(16,1).Range.Text = "74-56=18"
$t.Cell(16,2).Range.Text = "83-67=16"
$t.Cell(16,3).Range.Text = "90-41=49"
$t.Cell(16,4).Range.Text = "39+16=55"
$t.Cell(16,5).Range.Text = "25+68=93"

$t.Cell(17,1).Range.Text = "69+9=78"
$t.Cell(17,2).Range.Text = "12+49=61"
$t.Cell(17,3).Range.Text = "7+67=74"
$t.Cell(17,4).Range.Text = "91-83=8"
$t.Cell(17,5).Range.Text = "16+6=22"

$t.Cell(18,1).Range.Text = "41-17=24"
$t.Cell(18,2).Range.Text = "88-9=79"
$t.Cell(18,3).Range.Text = "61-43=18"
$t.Cell(18,4).Range.Text = "66-49=17"
$t.Cell(18,5).Range.Text = "69+3=72"

$t.Cell(19,1).Range.Text = "18+38=56"
$t.Cell(19,2).Range.Text = "81-13=68"
$t.Cell(19,3).Range.Text = "64+7=71"
$t.Cell(19,4).Range.Text = "15+7=22"
$t.Cell(19,5).Range.Text = "32-14=18"

$t.Cell(20,1).Range.Text = "58+8=66"
$t.Cell(20,2).Range.Text = "50-38=12"
$t.Cell(20,3).Range.Text = "9+53=62"
$t.Cell(20,4).Range.Text = "7+57=64"
$t.Cell(20,5).Range.Text = "13+69=82"
